# Deliverables Tracking.xlsx — "Updated Feedback, Architectural Deliverables,
# and Task Enumeration"
#
# 1. Fill in the Hours/Start/Due/Assigned-time columns for the three
#    existing "Architectural Deliverables" rows (Major Components BOM,
#    Hardware Block Diagram, Product Architecture).
# 2. Insert a new deliverable row ("Fill in Task Enumeration for Product
#    Definition ") right below them, with the same Hours/Start/Due/Assigned
#    pattern, pushing the remainder of the sheet down by one row.
# 3. Minor workbook-window / selection bookkeeping.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hardware Development Process")

# --- Fill in C/D/E/F for the three existing rows (11-13) ---------------
foreach ($r in 11..13) {
    $ws.Cells.Item($r, 3).Value = 2              # C: Estimate Task Hours
    $ws.Cells.Item($r, 4).Value = 42871           # D: Due Date
    $ws.Cells.Item($r, 4).NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 5).Value = 42878           # E: Class Assigned / done date
    $ws.Cells.Item($r, 5).NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 6).Value = "End of Day"    # F: Time
}

# --- Insert a new row for the "Task Enumeration for Product Definition" --
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = "Fill in Task Enumeration for Product Definition "
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 4).Value = 42871
$ws.Cells.Item(14, 4).NumberFormat = "d-mmm"
$ws.Cells.Item(14, 5).Value = 42878
$ws.Cells.Item(14, 5).NumberFormat = "d-mmm"
$ws.Cells.Item(14, 6).Value = "End of Day"

# --- Selection bookkeeping ----------------------------------------------
$ws.Range("D15").Select() | Out-Null
